$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the underlying data (B2:B4) ---
# The second and third runs now report the same timing as the first
# (40.2s) after adding new per-page context.
$ws.Range("B2").Value = 40.2
$ws.Range("B3").Value = 40.2
$ws.Range("B4").Value = 40.2

# --- Reposition / resize the embedded chart ---
$co = $ws.ChartObjects().Item(1)
$co.Top = 35.37952755905512
$co.Left = 227.6967716535433
$co.Width = 414.1365354330709
$co.Height = 255.26299212598425

# --- Update the saved selection/active cell ---
$ws.Range("B19").Select() | Out-Null
